$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 448
$ws.Range("F5").Value = 1339
$ws.Range("F6").Value = 7669
$ws.Range("F8").Value = 110
$ws.Range("F10").Value = 8462
$ws.Range("F13").Value = 65
$ws.Range("F14").Value = 5667
$ws.Range("F16").Value = 2625
$ws.Range("F17").Value = 1143
$ws.Range("F18").Value = 4595
$ws.Range("F23").Value = 535
$ws.Range("F24").Value = 3549
$ws.Range("F27").Value = 25
$ws.Range("F29").Value = 3046
$ws.Range("F30").Value = 40
$ws.Range("F31").Value = 111
$ws.Range("F32").Value = 346
$ws.Range("F33").Value = 128
$ws.Range("F34").Value = 314
$ws.Range("F35").Value = 571
$ws.Range("F36").Value = 660
$ws.Range("F39").Value = 1806
$ws.Range("F43").Value = 2944
$ws.Range("F45").Value = 2285
$ws.Range("F46").Value = 10
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 107
$ws.Range("F3").Value = 127
$ws.Range("F4").Value = 8
$ws.Range("F9").Value = 121
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1329
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1329
$ws.Range("F5").Value = 1339
$ws.Range("F6").Value = 7669
$ws.Range("F8").Value = 110
$ws.Range("F10").Value = 8462
$ws.Range("F12").Value = 65
$ws.Range("F13").Value = 5667
$ws.Range("F15").Value = 2625
$ws.Range("F16").Value = 1143
$ws.Range("F17").Value = 4595
$ws.Range("F20").Value = 107
$ws.Range("F22").Value = 127
$ws.Range("F23").Value = 535
$ws.Range("F24").Value = 8
$ws.Range("F25").Value = 3549
$ws.Range("F28").Value = 25
$ws.Range("F30").Value = 3046
$ws.Range("F31").Value = 346
$ws.Range("F32").Value = 128
$ws.Range("F33").Value = 314
$ws.Range("F35").Value = 571
$ws.Range("F36").Value = 660
$ws.Range("F40").Value = 1806
$ws.Range("F44").Value = 2944
$ws.Range("F45").Value = 2285
$ws.Range("F48").Value = 121
